$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are working with a client who has three distinct, yet related, tasks: sentiment analysis, topic classification, and spam detection from customer reviews. The dataset contains one million reviews, and the client's primary interest is in sentiment analysis.You are currently deciding on which multi-task learning architecture to use.What should you do?",
        "ques_type": 2,
        "options": [
            "Use a combination of hard parameter sharing for sentiment analysis and soft parameter sharing for topic classification and spam detection.",
            "Use hard parameter sharing for all tasks equally.",
            "Use separate deep learning models for each task.",
            "Use a shallow learning model due to the size of the dataset."
        ],
        "score": "Use a combination of hard parameter sharing for sentiment analysis and soft parameter sharing for topic classification and spam detection."
    },
    {
        "title": "Your company aims to build a multi-task learning model for a robotics system that guides a robot through a dynamic environment. You are currently trying to choose a strategy that handles potential overlap, redundancy, or conflicts between subtasks like object detection, path planning, collision avoidance, and energy management.What should you do?",
        "ques_type": 2,
        "options": [
            "Decompose the subtasks into primary and secondary subtasks on the basis of their importance.",
            "Combine all subtasks into a single complex task without decomposition.",
            "Treat all subtasks as separate without any correlation.",
            "Decompose the subtasks based solely on task difficulty, irrespective of their interaction."
        ],
        "score": "Decompose the subtasks into primary and secondary subtasks on the basis of their importance."
    },
    {
        "title": "You work for a startup that has a well-trained model for detecting tumors in x-ray images of lungs. You wish to extend its capabilities to detect fractures in x-rays of bones with a limited dataset, without compromising the tumor detection capabilities.What should you do?",
        "ques_type": 2,
        "options": [
            "Use the tumor detection model as a base and fine-tune it using the fracture dataset.",
            "Train a new model from scratch using the fracture dataset, then fine-tune it with the tumor dataset.",
            "Combine both datasets and retrain the model from scratch.",
            "Use an ensemble of two separate models: one for tumors and another for fractures."
        ],
        "score": "Use the tumor detection model as a base and fine-tune it using the fracture dataset."
    },
    {
        "title": "You are working on a robotic system designed for assisting with grocery store restocking. The model needs to concurrently handle tasks like identifying empty shelves, avoiding collisions with customers, and handling fragile items. You want to configure your loss function to prioritize customer safety.What should you do?",
        "ques_type": 2,
        "options": [
            "Increase the weight for collision avoidance loss while reducing weights for the other tasks.",
            "Use a dynamic weighting scheme that adjusts weights based on real-time store conditions.",
            "Introduce an adaptive feedback loop that increases collision avoidance loss weight if any near-miss incidents occur.",
            "Incorporate external safety metrics and datasets into the weighting scheme."
        ],
        "score": "Increase the weight for collision avoidance loss while reducing weights for the other tasks."
    }
]
'@

$ws.Range("A2").Delete() | Out-Null
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit() | Out-Null
